# Commit: "case with 380 kV done"
# Update the res_bus vm_pu results table (rows 2-25, cols B-N except G/H)
# with the voltage-magnitude values recomputed for the 380 kV case.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vmPuUpdates = @{
    2 = @{ "B"=1.02; "C"=1.056554295772895; "D"=1.061880159772376; "E"=1.063108978231971; "F"=1.074439240449316; "I"=1.054129719975071; "J"=1.061555325885056; "K"=1.064603020345378; "L"=1.065828508331149; "M"=1.077128450119488; "N"=1.024245067184439 }
    3 = @{ "B"=1.02; "C"=1.057537079297227; "D"=1.062652828369846; "E"=1.063968215666284; "F"=1.075342580909083; "I"=1.05444041915412; "J"=1.062190069750803; "K"=1.065190502055134; "L"=1.066502581242343; "M"=1.077848703985172; "N"=1.024460774060276 }
    4 = @{ "B"=1.02; "C"=1.058173449703919; "D"=1.063153164679208; "E"=1.064524916488953; "F"=1.07592785094246; "I"=1.054640423384994; "J"=1.062600595645402; "K"=1.065570347094107; "L"=1.066938822253564; "M"=1.078314865907837; "N"=1.024600184208667 }
    5 = @{ "B"=1.02; "C"=1.058441085194728; "D"=1.063363592931275; "E"=1.064759123839088; "F"=1.07617407629912; "I"=1.054724255887667; "J"=1.062773132957957; "K"=1.065729962429856; "L"=1.06712223391864; "M"=1.078510865509056; "N"=1.02465875195506 }
    6 = @{ "B"=1.02; "C"=1.058486028525383; "D"=1.063398929818458; "E"=1.064798458234848; "F"=1.076215429009851; "I"=1.054738317099918; "J"=1.062802099904611; "K"=1.065756758326581; "L"=1.067153030450272; "M"=1.078543776146524; "N"=1.024668583364029 }
    7 = @{ "B"=1.02; "C"=1.058177025449474; "D"=1.063155976092043; "E"=1.064528045311357; "F"=1.075931140320614; "I"=1.054641544539175; "J"=1.062602901286247; "K"=1.0655724801642; "L"=1.066941272947808; "M"=1.078317484766723; "N"=1.024600966952521 }
    8 = @{ "B"=1.02; "C"=1.056886339664937; "D"=1.062141209783629; "E"=1.063399212376903; "F"=1.074744372840765; "I"=1.054234936995178; "J"=1.061769880557947; "K"=1.064801623051684; "L"=1.066056299175529; "M"=1.077371839956455; "N"=1.024318000599235 }
    9 = @{ "B"=1.02; "C"=1.054615423081794; "D"=1.060355935807046; "E"=1.061415609849954; "F"=1.072658926758232; "I"=1.053510510370943; "J"=1.060300531111873; "K"=1.063441056108458; "L"=1.064497448293904; "M"=1.075706379713119; "N"=1.023818118793616 }
    10 = @{ "B"=1.02; "C"=1.053103839171154; "D"=1.059167758812218; "E"=1.060097011078111; "F"=1.071272599960863; "I"=1.053022258831806; "J"=1.059320035040105; "K"=1.062532573842687; "L"=1.063458666895702; "M"=1.074596733004101; "N"=1.023484040755719 }
    11 = @{ "B"=1.02; "C"=1.05244987654199; "D"=1.058653756001705; "E"=1.059526960708989; "F"=1.070673263057141; "I"=1.052809590024744; "J"=1.058895259601859; "K"=1.062138861087946; "L"=1.06300898250055; "M"=1.074116413513845; "N"=1.023339190703628 }
    12 = @{ "B"=1.02; "C"=1.052207051059386; "D"=1.058462906635977; "E"=1.059315356794976; "F"=1.07045078690467; "I"=1.05273040766208; "J"=1.058737447736646; "K"=1.061992569325535; "L"=1.062841967881194; "M"=1.073938027266085; "N"=1.023285358459553 }
    13 = @{ "B"=1.02; "C"=1.052259134075167; "D"=1.058503841097462; "E"=1.059360740280598; "F"=1.070498502246162; "I"=1.052747401033352; "J"=1.058771300323063; "K"=1.062023951623038; "L"=1.062877792245885; "M"=1.073976290537857; "N"=1.023296906941733 }
    14 = @{ "B"=1.02; "C"=1.052429802766293; "D"=1.058637978813743; "E"=1.059509466637684; "F"=1.070654870158785; "I"=1.052803048612033; "J"=1.058882215458704; "K"=1.06212676956629; "L"=1.062995176643561; "M"=1.074101667518742; "N"=1.023334741492494 }
    15 = @{ "B"=1.02; "C"=1.05253496871248; "D"=1.058720635351147; "E"=1.059601120200767; "F"=1.070751232744274; "I"=1.052837310064999; "J"=1.05895054980432; "K"=1.062190112585421; "L"=1.063067503483247; "M"=1.07417891986382; "N"=1.02335804884439 }
    16 = @{ "B"=1.02; "C"=1.053147252446009; "D"=1.059201881820641; "E"=1.060134862753915; "F"=1.071312396137876; "I"=1.0530363466105; "J"=1.059348221544802; "K"=1.062558696298811; "L"=1.063488513466932; "M"=1.074628613787938; "N"=1.0234936499617 }
    17 = @{ "B"=1.02; "C"=1.053531473158438; "D"=1.05950388588941; "E"=1.060469910276288; "F"=1.071664654814806; "I"=1.053160861954569; "J"=1.0595976137618; "K"=1.062789810322704; "L"=1.063752633305717; "M"=1.074910739768789; "N"=1.023578657801373 }
    18 = @{ "B"=1.02; "C"=1.0537556368879; "D"=1.059680086532886; "E"=1.060665425702099; "F"=1.071870213085971; "I"=1.053233368750032; "J"=1.059743059303135; "K"=1.062924582945008; "L"=1.063906700953091; "M"=1.07507531491745; "N"=1.023628222880533 }
    19 = @{ "B"=1.02; "C"=1.053832080112431; "D"=1.059740174297089; "E"=1.060732106277774; "F"=1.071940318667212; "I"=1.053258071171939; "J"=1.059792648934725; "K"=1.06297053142553; "L"=1.063959235878999; "M"=1.075131433445123; "N"=1.023645120133092 }
    20 = @{ "B"=1.02; "C"=1.053490244291497; "D"=1.059471478863655; "E"=1.060433953749904; "F"=1.071626851293739; "I"=1.053147515137978; "J"=1.059570858474287; "K"=1.062765017307614; "L"=1.063724294600588; "M"=1.074880468675482; "N"=1.023569539181291 }
    21 = @{ "B"=1.02; "C"=1.052379542728326; "D"=1.05859847654268; "E"=1.059465666608356; "F"=1.070608819713426; "I"=1.052786666964457; "J"=1.058849554574663; "K"=1.062096493599137; "L"=1.062960609354569; "M"=1.074064746402409; "N"=1.023323600943663 }
    22 = @{ "B"=1.02; "C"=1.051681694785616; "D"=1.058050015017979; "E"=1.058857666090035; "F"=1.069969578121301; "I"=1.05255870167937; "J"=1.058395861186951; "K"=1.061675881312107; "L"=1.06248055582162; "M"=1.073552019198617; "N"=1.023168805188422 }
    23 = @{ "B"=1.02; "C"=1.052051590087969; "D"=1.058340723720277; "E"=1.059179902427459; "F"=1.070308372586042; "I"=1.052679653135632; "J"=1.058636389517226; "K"=1.061898882562624; "L"=1.062735030949939; "M"=1.073823811050306; "N"=1.023250880828186 }
    24 = @{ "B"=1.02; "C"=1.053508873674067; "D"=1.05948612205922; "E"=1.060450200685041; "F"=1.071643932794317; "I"=1.053153546363473; "J"=1.059582948101123; "K"=1.062776220303729; "L"=1.063737099607655; "M"=1.074894146830695; "N"=1.023573659550565 }
    25 = @{ "B"=1.02; "C"=1.055202097070765; "D"=1.060817123787615; "E"=1.061927753822099; "F"=1.073197370417134; "I"=1.05369872866314; "J"=1.060680560871121; "K"=1.063793052612299; "L"=1.064900373144459; "M"=1.076136829442669; "N"=1.023947497050901 }
}

foreach ($rowKey in $vmPuUpdates.Keys) {
    $rowValues = $vmPuUpdates[$rowKey]
    foreach ($colKey in $rowValues.Keys) {
        $cellRef = "{0}{1}" -f $colKey, $rowKey
        $ws.Range($cellRef).Value = $rowValues[$colKey]
    }
}
